# infrastructure.xlsx: add Monitor/Printer sheets, extend Device sheet
# with new columns, and drop now-redundant columns from Computer/Harddisk.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Device sheet: insert a "type" column right after "name", and append
# eight new trailing columns (location, ip, mac, website, driver,
# support, comment, computer).
# ---------------------------------------------------------------------
$device = $wb.Worksheets.Item("Device")
$device.Range("B:B").Insert()
$device.Cells.Item(1, 2).Value = "type"
$device.Cells.Item(1, 10).Value = "location"
$device.Cells.Item(1, 11).Value = "ip"
$device.Cells.Item(1, 12).Value = "mac"
$device.Cells.Item(1, 13).Value = "website"
$device.Cells.Item(1, 14).Value = "driver"
$device.Cells.Item(1, 15).Value = "support"
$device.Cells.Item(1, 16).Value = "comment"
$device.Cells.Item(1, 17).Value = "computer"

# ---------------------------------------------------------------------
# Computer sheet: drop the "location" column (now tracked per-device).
# ---------------------------------------------------------------------
$computer = $wb.Worksheets.Item("Computer")
$computer.Range("C:C").Delete()

# ---------------------------------------------------------------------
# Harddisk sheet: drop the "computer" column (now tracked on Device).
# ---------------------------------------------------------------------
$harddisk = $wb.Worksheets.Item("Harddisk")
$harddisk.Range("D:D").Delete()

# ---------------------------------------------------------------------
# New Monitor sheet, appended after Harddisk.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$monitor = $wb.Worksheets.Add($null, $lastSheet)
$monitor.Name = "Monitor"
$monitor.Cells.Item(1, 1).Value = "resolution"

# ---------------------------------------------------------------------
# New Printer sheet, appended after Monitor (left empty).
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$printer = $wb.Worksheets.Add($null, $lastSheet)
$printer.Name = "Printer"
